$d = $word.ActiveDocument

# Locate the paragraph that contains the ellipsis ("...") placeholder text -
# this is the paragraph the author appended their "M name is Tanvir.Studying
# SCU Sydney" sentence to.
$search = $d.Content
$found = $search.Find.Execute("…", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target paragraph (ellipsis run) to edit."
}

$target = $search.Paragraphs(1).Range

# Rebuild the whole paragraph as OOXML so we can add the new runs as their
# own <w:r> elements (split around a spell-check proofErr pair for the
# run-on word "Tanvir.Studying"), exactly like Word would when a user types
# new text at the end of an existing paragraph and the proofer flags one of
# the new words.
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
    'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
    'w14:paraId="69584E20" w14:textId="5B6E673B" w:rsidR="00EF6287" w:rsidRDefault="00EF6287">' +
    '<w:r><w:t>…</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">M name is </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Tanvir.Studying</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> SCU Sydney</w:t></w:r>' +
    '</w:p>'

$target.InsertXML($newParaXml)
